# "pumped storage and bio-mass"
# - Fix the swapped D1/E1 header labels (D1 = Year, E1 = Capacity (MW))
# - Insert a new "Greensville" natural-gas plant row (row 31), pushing the
#   remaining Natural gas / Combined Cycle / Oil / Uranium / Waste rows down
# - Backfill "Year" (in-service / commission year) values in column F for
#   several plants that previously had no year/footnote entry
# - Drop the stray "[1]" footnote that used to sit next to Aria Energy

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header row: D should read "Year", E should read "Capacity (MW)" ---
$ws.Range("D1").Value = "Year"
$ws.Range("E1").Value = "Capacity (MW)"

# --- Insert the new Greensville row at position 31 (shifts 31-40 -> 32-41) ---
$ws.Rows(31).Insert()

$ws.Range("A31").Value = "Greensville"
$ws.Range("C31").Value = "Natural gas"
$ws.Range("E31").Value = 1588
$ws.Range("F31").Value = 2018

# --- Remove the stray "[1]" footnote on the Aria Energy row ---
$ws.Range("F19").ClearContents()

# --- Backfill commission/in-service years in column F ---
$ws.Range("F12").Value = 1992    # Chesterfield Power Station
$ws.Range("F20").Value = 2011    # Bear Garden
$ws.Range("F21").Value = 2016    # Brunswick County
$ws.Range("F22").Value = 1990    # Darbytown
$ws.Range("F24").Value = 1992    # Elizabeth River
$ws.Range("F25").Value = 1994    # Gordonsville
$ws.Range("F27").Value = 2001    # Ladysmith
$ws.Range("F29").Value = 2000    # Remington
$ws.Range("F32").Value = 2014    # Warren County
$ws.Range("F36").Value = 2003    # Possum Point Power Station

# --- Restore the selection state (C16:C17 on the Bath County / Smith
#     Mountain pumped-storage rows) ---
$ws.Range("C16:C17").Select()
